# Why Scientific Notation.pptx - date placeholder refresh + answer tweak
#
# 1. The "datetimeFigureOut" Date placeholder on the Slide Master and on
#    every Custom Layout gets its cached text bumped from 2/10/2019 to
#    2/5/2020 (this happens automatically in real PowerPoint whenever the
#    deck is opened/saved on a later date with the field set to auto-update).
# 2. The big "32.5" answer shown on slides 1 and 2 becomes "32.4".

$p = $ppt.ActivePresentation

$oldDate = "2/10/2019"
$newDate = "2/5/2020"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if (-not $isDate) { continue }

        $tr = $shp.TextFrame.TextRange
        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi)
            if ($para.Text -eq $oldDate) {
                $para.Text = $newDate
            }
        }
    }
}

# Slide Master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Custom Layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide content: "32.5" -> "32.4" (slides 1 and 2).
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi)
            for ($ri = 1; $ri -le $para.Runs().Count; $ri++) {
                $run = $para.Runs($ri)
                if ($run.Text -eq "32.5") {
                    $run.Text = "32.4"
                }
            }
        }
    }
}

Write-Output "edit complete"
